$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update GoalPoseX (B12) and GoalPoseY (B13) values
$ws.Range("B12").Value = 154
$ws.Range("B13").Value = -339

# Move the active cell selection to F31
$ws.Range("F31").Select()
